$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.723.25"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.599.89"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'211.76"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.824.17"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.604.22"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'208.51"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "'143.62"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "'15.33"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'1.28"
$ws.Range("E33").Value = "  +19.68%  "
$ws.Range("D34").Value = "1.279.61"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "'2.49"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'0.591"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "'62.69"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "1.735.85"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'90.41"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'7.54"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +1.53%  "
